$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the last existing table row (row 64) down onto the
# new row (65) so the new cells inherit the same number format / hyperlink
# font styling as the rest of the table, then overwrite with the new values.
$ws.Range("B64:F64").Copy()
$ws.Range("B65:F65").PasteSpecial(-4122)

$ws.Range("B65").Value = 55
$ws.Range("E65").Value = "https://programmingport.hashnode.dev/case-statement-or-shell-scripting"
$ws.Range("F65").Value = "https://dev.to/rahulmishra05/case-statement-shell-scripting-2o2a"
$ws.Range("C65").Value = "Case Statement | Shell Scripting"
$ws.Range("D65").Value = "12/12/2020"

# Grow the table (ListObject) so the new row becomes part of Table2 / its
# AutoFilter range, matching B10:F65.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B10:F65"))

# Update the view: scroll position moved right one column and selection
# moved to the cell just past the newly-added row.
$ws.Range("G65").Select()
